# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-03 to 2023-10-04 (Excel serial date 45202 -> 45203).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

$ws.Range("C2:C$lastRow").Value = 45203
